$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-23: column A previously held bogus/misplaced values (parsed as dates or shifted
# text). Replace them with the correct label text that belongs on each of these rows.
$ws.Range("A2").Value = '                National Gem and Jewellery Authority'
$ws.Range("A3").Value = '               Sri Lanka Customs'
$ws.Range("A4").Value = '              Central Bank of Sri Lanka'
$ws.Range("A5").Value = '(a)  The latest version of SITC Revision 4 published in 2006'
$ws.Range("A6").Value = '(b) Provisional'
$ws.Range("A7").Value = '2.04: Import Performance based on Standard International Trade Classification (SITC) Monthly 2014-2024 (a)'
$ws.Range("A8").Value = '2.04: Import Performance based on Standard International Trade Classification (SITC) Monthly 2014-2024 (a)'
$ws.Range("A9").Value = '2023'
$ws.Range("A10").Value = '2024'
$ws.Range("A11").Value = 'January'
$ws.Range("A12").Value = 'January'
$ws.Range("A13").Value = 'January'
$ws.Range("A14").Value = 'January'
$ws.Range("A15").Value = 'January'
$ws.Range("A16").Value = 'January'
$ws.Range("A17").Value = 'January'
$ws.Range("A18").Value = 'January'
$ws.Range("A19").Value = 'January'
$ws.Range("A20").Value = 'January'
$ws.Range("A21").Value = 'January'
$ws.Range("A22").Value = 'Sources: Ceylon Petroleum Corporation and Other Exporters of Petroleum'
$ws.Range("A23").Value = 'Table 2.04.3: Imports (US$ Million)'

# Rows 24-161: column A held leftover/duplicated text that does not belong here; clear it.
$ws.Range("A24:A161").ClearContents()

# Rows 162-323 were entirely bogus trailing rows; remove them so the sheet ends at row 161.
$ws.Range("A162:A323").EntireRow.Delete()
